$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Trf", "Tfr2", "FAPs", 3, 1, 1.977922333333334, 5.933767, 0.01156576927021167, 0.01156576927021167, 3, 1, 1.407309666666666, 4.221928999999999, 0.7536752551149739, 0.7536752551149739, 2.783549219615889, 25.051942976543, 0.008716834105327706, 0.008716834105327705)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }
$row3 = @("ECs", "Trf", "Tfr2", "Inflammatory-Mac", 3, 1, 1.977922333333334, 5.933767, 0.01156576927021167, 0.01156576927021167, 1, 0.3333333333333333, 0.1019666666666667, 0.3059, 0.05460756458473616, 0.05460756458473616, 0.2016821472555556, 1.8151393253, 0.0006315784923952406, 0.0006315784923952405)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }
$row4 = @("ECs", "Trf", "Tfr2", "MuSCs", 3, 1, 1.977922333333334, 5.933767, 0.01156576927021167, 0.01156576927021167, 3, 1, 0.3265433333333334, 0.97963, 0.1748780924947535, 0.1748780924947535, 0.6458773518011113, 5.812896166210001, 0.002022599668209054, 0.002022599668209054)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }
$row5 = @("ECs", "Trf", "Tfr2", "Resolving-Mac", 3, 1, 1.977922333333334, 5.933767, 0.01156576927021167, 0.01156576927021167, 1, 0.3333333333333333, 0.031443, 0.094329, 0.01683908780553638, 0.01683908780553637, 0.062191811927, 0.559726307343, 0.0001947570042796687, 0.0001947570042796686)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }
$row6 = @("FAPs", "Trf", "Tfr2", "FAPs", 3, 1, 14.16077566666667, 42.482327, 0.08280419371769795, 0.08280419371769795, 3, 1, 1.407309666666666, 4.221928999999999, 0.7536752551149739, 0.7536752551149739, 19.92859648319811, 179.357368348783, 0.06240747182477573, 0.06240747182477573)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row6[$i] }
$row7 = @("FAPs", "Trf", "Tfr2", "Inflammatory-Mac", 3, 1, 14.16077566666667, 42.482327, 0.08280419371769795, 0.08280419371769795, 1, 0.3333333333333333, 0.1019666666666667, 0.3059, 0.05460756458473616, 0.05460756458473616, 1.443927092144444, 12.9953438293, 0.004521735356326195, 0.004521735356326195)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row7[$i] }
$row8 = @("FAPs", "Trf", "Tfr2", "MuSCs", 3, 1, 14.16077566666667, 42.482327, 0.08280419371769795, 0.08280419371769795, 3, 1, 0.3265433333333334, 0.97963, 0.1748780924947535, 0.1748780924947535, 4.624106888778889, 41.61696199900999, 0.01448063944791707, 0.01448063944791707)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row8[$i] }
$row9 = @("FAPs", "Trf", "Tfr2", "Resolving-Mac", 3, 1, 14.16077566666667, 42.482327, 0.08280419371769795, 0.08280419371769795, 1, 0.3333333333333333, 0.031443, 0.094329, 0.01683908780553638, 0.01683908780553637, 0.4452572692869999, 4.007315423583, 0.001394347088678959, 0.001394347088678959)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row9[$i] }
$row10 = @("Inflammatory-Mac", "Trf", "Tfr2", "FAPs", 3, 1, 80.27561433333334, 240.826843, 0.4694063147763454, 0.4694063147763453, 3, 1, 1.407309666666666, 4.221928999999999, 0.7536752551149739, 0.7536752551149739, 112.9726480489052, 1016.753832440147, 0.3537799240416418, 0.3537799240416418)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row10[$i] }
$row11 = @("Inflammatory-Mac", "Trf", "Tfr2", "Inflammatory-Mac", 3, 1, 80.27561433333334, 240.826843, 0.4694063147763454, 0.4694063147763453, 1, 0.3333333333333333, 0.1019666666666667, 0.3059, 0.05460756458473616, 0.05460756458473616, 8.18543680818889, 73.6689312737, 0.02563313565063227, 0.02563313565063226)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row11[$i] }
$row12 = @("Inflammatory-Mac", "Trf", "Tfr2", "MuSCs", 3, 1, 80.27561433333334, 240.826843, 0.4694063147763454, 0.4694063147763453, 3, 1, 0.3265433333333334, 0.97963, 0.1748780924947535, 0.1748780924947535, 26.21346668978778, 235.92120020809, 0.08208888093307909, 0.08208888093307908)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row12[$i] }
$row13 = @("Inflammatory-Mac", "Trf", "Tfr2", "Resolving-Mac", 3, 1, 80.27561433333334, 240.826843, 0.4694063147763454, 0.4694063147763453, 1, 0.3333333333333333, 0.031443, 0.094329, 0.01683908780553638, 0.01683908780553637, 2.524106141483, 22.716955273347, 0.007904374150992127, 0.007904374150992124)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row13[$i] }
$row14 = @("MuSCs", "Trf", "Tfr2", "FAPs", 3, 1, 0.3725693333333333, 1.117708, 0.002178574392872141, 0.002178574392872141, 3, 1, 1.407309666666666, 4.221928999999999, 0.7536752551149739, 0.7536752551149739, 0.5243204243035554, 4.718883818731999, 0.00164193761133486, 0.00164193761133486)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row14[$i] }
$row15 = @("MuSCs", "Trf", "Tfr2", "Inflammatory-Mac", 3, 1, 0.3725693333333333, 1.117708, 0.002178574392872141, 0.002178574392872141, 1, 0.3333333333333333, 0.1019666666666667, 0.3059, 0.05460756458473616, 0.05460756458473616, 0.03798965302222222, 0.3419068772, 0.0001189666418614178, 0.0001189666418614178)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row15[$i] }
$row16 = @("MuSCs", "Trf", "Tfr2", "MuSCs", 3, 1, 0.3725693333333333, 1.117708, 0.002178574392872141, 0.002178574392872141, 3, 1, 0.3265433333333334, 0.97963, 0.1748780924947535, 0.1748780924947535, 0.1216600320044444, 1.09494028804, 0.0003809849341833956, 0.0003809849341833956)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row16[$i] }
$row17 = @("MuSCs", "Trf", "Tfr2", "Resolving-Mac", 3, 1, 0.3725693333333333, 1.117708, 0.002178574392872141, 0.002178574392872141, 1, 0.3333333333333333, 0.031443, 0.094329, 0.01683908780553638, 0.01683908780553637, 0.011714697548, 0.105432277932, [double]"3.668520549246708e-05", [double]"3.668520549246707e-05")
for ($i = 0; $i -lt $row17.Length; $i++) { $ws.Cells.Item(17, $i + 1).Value = $row17[$i] }
$row18 = @("Resolving-Mac", "Trf", "Tfr2", "FAPs", 3, 1, 74.22831733333334, 222.684952, 0.434045147842873, 0.4340451478428729, 3, 1, 1.407309666666666, 4.221928999999999, 0.7536752551149739, 0.7536752551149739, 104.4622285236009, 940.1600567124079, 0.3271290875318938, 0.3271290875318938)
for ($i = 0; $i -lt $row18.Length; $i++) { $ws.Cells.Item(18, $i + 1).Value = $row18[$i] }
$row19 = @("Resolving-Mac", "Trf", "Tfr2", "Inflammatory-Mac", 3, 1, 74.22831733333334, 222.684952, 0.434045147842873, 0.4340451478428729, 1, 0.3333333333333333, 0.1019666666666667, 0.3059, 0.05460756458473616, 0.05460756458473616, 7.568814090755556, 68.1193268168, 0.02370214844352104, 0.02370214844352104)
for ($i = 0; $i -lt $row19.Length; $i++) { $ws.Cells.Item(19, $i + 1).Value = $row19[$i] }
$row20 = @("Resolving-Mac", "Trf", "Tfr2", "MuSCs", 3, 1, 74.22831733333334, 222.684952, 0.434045147842873, 0.4340451478428729, 3, 1, 0.3265433333333334, 0.97963, 0.1748780924947535, 0.1748780924947535, 24.23876216975111, 218.14885952776, 0.07590498751136489, 0.07590498751136487)
for ($i = 0; $i -lt $row20.Length; $i++) { $ws.Cells.Item(20, $i + 1).Value = $row20[$i] }
$row21 = @("Resolving-Mac", "Trf", "Tfr2", "Resolving-Mac", 3, 1, 74.22831733333334, 222.684952, 0.434045147842873, 0.4340451478428729, 1, 0.3333333333333333, 0.031443, 0.094329, 0.01683908780553638, 0.01683908780553637, 2.333960981912, 21.005648837208, 0.007308924356093155, 0.007308924356093153)
for ($i = 0; $i -lt $row21.Length; $i++) { $ws.Cells.Item(21, $i + 1).Value = $row21[$i] }
